$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) cells that change to remain stored as text,
# matching the source data which keeps prices like "1.00" / "0.0000271"
# as literal strings (not numbers where trailing zeros would be lost).
$ws.Range('D2:D7').NumberFormat = '@'
$ws.Range('D9:D12').NumberFormat = '@'
$ws.Range('D14:D51').NumberFormat = '@'

$ws.Range('D2').Value = '81.445.70'
$ws.Range('E2').Value = '  +5.85%  '
$ws.Range('D3').Value = '3.211.22'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '211.44'
$ws.Range('E5').Value = '  +5.15%  '
$ws.Range('D6').Value = '639.36'
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').Value = '0.298'
$ws.Range('E7').Value = '  +31.26%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +4.29%  '
$ws.Range('D10').Value = '3.203.19'
$ws.Range('E10').Value = '  +2.15%  '
$ws.Range('D11').Value = '0.597'
$ws.Range('E11').Value = '  +14.21%  '
$ws.Range('D12').Value = '0.0000271'
$ws.Range('E12').Value = '  +21.57%  '
$ws.Range('E13').Value = '  +2.25%  '
$ws.Range('D14').Value = '5.38'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Value = '3.798.05'
$ws.Range('E15').Value = '  +2.40%  '
$ws.Range('D16').Value = '32.32'
$ws.Range('E16').Value = '  +5.93%  '
$ws.Range('D17').Value = '81.143.38'
$ws.Range('E17').Value = '  +5.74%  '
$ws.Range('D18').Value = '3.208.84'
$ws.Range('E18').Value = '  +3.08%  '
$ws.Range('D19').Value = '14.48'
$ws.Range('E19').Value = '  +3.56%  '
$ws.Range('D20').Value = '3.17'
$ws.Range('E20').Value = '  +13.17%  '
$ws.Range('D21').Value = '449.30'
$ws.Range('E21').Value = '  +12.21%  '
$ws.Range('D22').Value = '9.32'
$ws.Range('E22').Value = '  +1.23%  '
$ws.Range('D23').Value = '5.27'
$ws.Range('E23').Value = '  +12.52%  '
$ws.Range('D24').Value = '7.15'
$ws.Range('E24').Value = '  +6.40%  '
$ws.Range('D25').Value = '5.06'
$ws.Range('E25').Value = '  +10.12%  '
$ws.Range('D26').Value = '3.368.80'
$ws.Range('E26').Value = '  +2.30%  '
$ws.Range('D27').Value = '77.64'
$ws.Range('E27').Value = '  +4.42%  '
$ws.Range('D28').Value = '11.22'
$ws.Range('E28').Value = '  +7.86%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0000128'
$ws.Range('E29').Value = '  +12.00%  '
$ws.Range('B30').Value = 'Dai'
$ws.Range('C30').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '9.24'
$ws.Range('E31').Value = '  +7.15%  '
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').Value = '565.50'
$ws.Range('E33').Value = '  +8.82%  '
$ws.Range('D34').Value = '1.50'
$ws.Range('E34').Value = '  +2.67%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').Value = '2.05'
$ws.Range('E35').Value = '  +5.34%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = '0.153'
$ws.Range('E36').Value = '  +13.69%  '
$ws.Range('D37').Value = '23.19'
$ws.Range('E37').Value = '  +7.38%  '
$ws.Range('D38').Value = '0.126'
$ws.Range('E38').Value = '  +22.28%  '
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.04%  '
$ws.Range('D40').Value = '0.415'
$ws.Range('E40').Value = '  +7.06%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').Value = '2.09'
$ws.Range('E41').Value = '  +20.25%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = '20.82'
$ws.Range('E42').Value = '  +3.80%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').Value = '5.86'
$ws.Range('E43').Value = '  +10.25%  '
$ws.Range('B44').Value = 'Monero'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D44').Value = '158.57'
$ws.Range('E44').Value = '  -3.26%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '191.01'
$ws.Range('E46').Value = '  -2.22%  '
$ws.Range('D47').Value = '2.89'
$ws.Range('E47').Value = '  +16.18%  '
$ws.Range('D48').Value = '1.35'
$ws.Range('E48').Value = '  +5.43%  '
$ws.Range('D49').Value = '0.792'
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '4.33'
$ws.Range('E50').Value = '  +6.55%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '43.03'
$ws.Range('E51').Value = '  +1.29%  '
